# Script 1 - atualização automática de dados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data (Ano, BR Deflacionado, NE Deflacionado, SE Deflacionado)
$data = @(
    @(2010, 23257.86610669455, 11418.82435341108, 13560.43311809023),
    @(2011, 25192.59314690959, 12245.94221027538, 14257.77027193148),
    @(2012, 26979.84043437291, 13211.21379819802, 15697.30052998581),
    @(2013, 27976.04665397691, 13739.79027459498, 16072.90969847004),
    @(2014, 29912.95664259601, 14745.8907374447,  16785.95219769734),
    @(2015, 31911.5976261855,  15975.15151436092, 17673.68733391674),
    @(2016, 34224.49724617274, 17607.09315882314, 18604.68706838786),
    @(2017, 35211.13679492509, 18274.86689908277, 19513.7812222432),
    @(2018, 36646.18938775777, 19083.90609078256, 20594.91845447646),
    @(2019, 37893.92756225695, 19561.56719637316, 20960.20827855906),
    @(2020, 40040.12552822272, 20906.9003924803,  21326.36830874367),
    @(2021, 44932.83670743866, 22962.50837336785, 23155.15872464016),
    @(2022, 48783.70604975082, 24759.83037637985, 25299.05850862464),
    @(2023, 51300.70579350938, 26237.41536180414, 26006.98661973922)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value2 = $entry[0]
    $ws.Cells.Item($row, 2).Value2 = $entry[1]
    $ws.Cells.Item($row, 3).Value2 = $entry[2]
    $ws.Cells.Item($row, 4).Value2 = $entry[3]
    $row = $row + 1
}
